# 赵行保加班调休统计.xlsx — add the August ("八月加班调休") sheet and
# correct two mistyped time entries on the July ("七月加班调休") sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fix the two "9:00-10:00" -> "9:00-11:00" typos on 七月加班调休 and
#    bump the matching duration cells (2h instead of 1h). The monthly
#    totals / carry-over formulas recalc automatically.
# ---------------------------------------------------------------------
$julySheet = $wb.Worksheets.Item("七月加班调休")
$julySheet.Range("E5").Value = "9:00-11:00"
$julySheet.Range("F5").Value = 2
$julySheet.Range("E7").Value = "9:00-11:00"
$julySheet.Range("F7").Value = 2

# ---------------------------------------------------------------------
# 2) Add the new 八月加班调休 sheet, right after 七月加班调休, by copying
#    the July sheet so all formatting / merged cells / number formats
#    carry over identically, then clearing the old data and entering
#    August's entries.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$julySheet.Copy($null, $wb.Worksheets.Item($sheetCount))
$augSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$augSheet.Name = "八月加班调休"

# Wipe the copied July data rows (keeps the per-cell styling intact).
$augSheet.Range("A3:F18").ClearContents()

# First (only, so far) August entry.
$augSheet.Range("A3").Value = 43678
$augSheet.Range("B3").Value = "19:00-00:30"
$augSheet.Range("C3").Value = 5.5
$augSheet.Range("D3").Value = 43679
$augSheet.Range("E3").Value = "9:00-12:00"
$augSheet.Range("F3").Value = 3

# Carry-over total formula references July's new closing balance (39).
$augSheet.Range("F21").Formula = "=39+F20"

# Cosmetic formatting to mirror the new sheet's column layout / margins.
$augSheet.Columns.Item(2).ColumnWidth = 12.25
$augSheet.Columns.Item(3).ColumnWidth = 10.875
$augSheet.Columns.Item(5).ColumnWidth = 11.25
$augSheet.Columns.Item(6).ColumnWidth = 10.875
$augSheet.PageSetup.LeftMargin = 54
$augSheet.PageSetup.RightMargin = 54
$augSheet.PageSetup.TopMargin = 72
$augSheet.PageSetup.BottomMargin = 72
$augSheet.PageSetup.HeaderMargin = 36
$augSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 3) Selection bookkeeping: July is no longer the active tab/selection,
#    August is, with its own last-used cell.
# ---------------------------------------------------------------------
$julySheet.Range("I17").Select()
$augSheet.Range("I19").Select()
$augSheet.Activate()
